$d = $word.ActiveDocument

# Locate the "Simple Sequence Diagram" heading text (the real Heading 1
# paragraph/bookmark - the Table of Contents entry that shows the same
# words is a separate field result and is not matched by Find).
$rng = $d.Content
$found = $rng.Find.Execute("Simple Sequence Diagram", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $headStart = $rng.Start
    $newWord = "Example"

    # Replace the whole heading range with the new first word. Replacing
    # across the full original range (rather than just the "Simple" word)
    # keeps the paragraph's bookmark anchored before the new text instead
    # of being pushed in between the edited pieces.
    $rng.Text = $newWord

    # Re-insert the untouched remainder (" Sequence Diagram") right after
    # the word we just replaced; this creates it as its own run, matching
    # how Word splits a partially-edited run into separate pieces.
    $afterNewWord = $headStart + $newWord.Length
    $tail = $d.Range($afterNewWord, $afterNewWord)
    $tail.InsertAfter(" Sequence Diagram")
}
